$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Fact" + "s" + ":" (three separate runs) -> a single "Facts:" run.
#    A same-text Find/Replace over that exact paragraph forces Word to
#    re-write the paragraph's runs as one consolidated run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Facts:", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Facts:", 2)

# ---------------------------------------------------------------------------
# 2) Drop the stray "_GoBack" bookmark from its old spot (the empty
#    paragraph right after the "...pests control." bullet).  It gets
#    re-created further down, at the end of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) After the "From here...price and temperature relationship." paragraph,
#    add: a blank paragraph, the new "database use" paragraph, and a
#    paragraph that (re)houses the "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*From here*") {
        $anchorPara = $d.Paragraphs.Item($i)
    }
}

$r = $anchorPara.Range.Duplicate
$r.Collapse(0)
$null = $r.InsertParagraphAfter()

$blankIndex = $anchorPara.Index + 1
$r2 = $d.Paragraphs.Item($blankIndex).Range.Duplicate
$r2.Collapse(0)
$null = $r2.InsertParagraphAfter()

$textIndex = $blankIndex + 1
$d.Paragraphs.Item($textIndex).Range.Text = "One use for our database could be to predict what effect current weather conditions should have on the price of potato futures. If temperatures were outside of the optimal range for potato growth for x number of days, we should be able to make an accurate forecast of the future prices of potatoes for that season. This information would allow us to determine if the future prices of potatoes are high or low and we would be able to buy and sell futures accordingly. We have eight years of production, pricing and weather data to analyze and use for our forecasting. "

$r3 = $d.Paragraphs.Item($textIndex).Range.Duplicate
$r3.Collapse(0)
$null = $r3.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4) Re-create "_GoBack" inside that new, otherwise-empty paragraph.
#    A zero-length bookmark dropped exactly on a freshly split paragraph
#    boundary tends to bleed into the following paragraph, so: type a
#    placeholder character, bookmark it, then erase the character - the
#    bookmark collapses back down but stays anchored inside this paragraph.
# ---------------------------------------------------------------------------
$bmIndex = $textIndex + 1
$bmPara = $d.Paragraphs.Item($bmIndex)
$bmPara.Range.Text = "X"
$bmRange = $bmPara.Range.Duplicate
$null = $bmRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$clearRange = $d.Bookmarks.Item("_GoBack").Range.Duplicate
$clearRange.Text = ""
